# Adapt column header formatting to respective input file names.
# Headers suffixed "_old"/"_new" become "_FV2404"/"_FV2410"; "diff" stays.
# Also: wrap the data range A1:U60 in an Excel Table ("Table1") and freeze
# the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Create the table on a small, unformatted scratch range first, then
#    resize it onto the real A1:U60 range. Doing it this way (rather than
#    adding the table directly on top of the already-bold header row)
#    avoids Excel baking the header row's existing bold formatting into a
#    new table "headerRowDxf" / dxfs entry - the real header cells keep
#    their original style untouched.
# ---------------------------------------------------------------------
$scratch = $ws.Range("AA1:AA2")
$scratch.Cells.Item(1, 1).Value = "scratch1"
$scratch.Cells.Item(2, 1).Value = "scratch2"

$tbl = $ws.ListObjects.Add(1, $scratch, [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"

$dataRange = $ws.Range("A1:U60")
$tbl.Resize($dataRange)

# remove the scratch values again, they are no longer needed
$ws.Range("AA1:AA2").ClearContents()

# ---------------------------------------------------------------------
# 2) Rename the header cells. Writing straight into the header row cells
#    (which belong to the table now) keeps both the worksheet cell text
#    and the table's column names in sync.
# ---------------------------------------------------------------------
$headerMap = @{
    "Segmentname_old"          = "Segmentname_FV2404"
    "Segmentgruppe_old"        = "Segmentgruppe_FV2404"
    "Segment_old"              = "Segment_FV2404"
    "Datenelement_old"         = "Datenelement_FV2404"
    "Segment ID_old"           = "Segment ID_FV2404"
    "Code_old"                 = "Code_FV2404"
    "Qualifier_old"            = "Qualifier_FV2404"
    "Beschreibung_old"         = "Beschreibung_FV2404"
    "Bedingungsausdruck_old"   = "Bedingungsausdruck_FV2404"
    "Bedingung_old"            = "Bedingung_FV2404"
    "diff"                     = "diff"
    "Segmentname_new"          = "Segmentname_FV2410"
    "Segmentgruppe_new"        = "Segmentgruppe_FV2410"
    "Segment_new"              = "Segment_FV2410"
    "Datenelement_new"         = "Datenelement_FV2410"
    "Segment ID_new"           = "Segment ID_FV2410"
    "Code_new"                 = "Code_FV2410"
    "Qualifier_new"            = "Qualifier_FV2410"
    "Beschreibung_new"         = "Beschreibung_FV2410"
    "Bedingungsausdruck_new"   = "Bedingungsausdruck_FV2410"
    "Bedingung_new"            = "Bedingung_FV2410"
}

for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $current = [string]$cell.Value2
    if ($headerMap.ContainsKey($current)) {
        $cell.Value = $headerMap[$current]
    }
}

# ---------------------------------------------------------------------
# 3) Freeze the header row (split below row 1).
# ---------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
